# Refresh the cryptos list (coinranking.com snapshot) in place.
# Updates Price (D) and Volume(1h) (E) for every coin row, and for the
# rows whose ranking position swapped, the Coin (B) and Link (C) too.
# D/E columns are stored as plain text (e.g. "3.100.84", "  -6.28%  "),
# so force a text NumberFormat before writing to stop Excel COM from
# reinterpreting the strings as numbers/percentages and losing
# formatting (trailing zeros, thousand-dot grouping, padding).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '56.980.15'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -5.85%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.107.39'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -6.18%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '518.29'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -7.42%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '132.20'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -7.89%  '

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.05%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.099.44'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -6.45%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.442'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -7.22%  '

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -10.21%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.107'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -10.32%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.378'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -8.15%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.647.47'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -5.73%  '

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.40%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.25'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -7.29%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.117.57'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -5.34%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '57.016.72'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -5.60%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0000148'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -11.16%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.70'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -7.59%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.80'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -11.36%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.88'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -8.82%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '342.89'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -8.86%  '

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.06%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '67.90'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -8.21%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.499'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -8.79%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.262.43'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -5.07%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.75%  '

$ws.Range('B28').Value = 'PEPE'
$ws.Range('C28').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0925'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -11.38%  '

$ws.Range('B29').Value = 'Kaspa'
$ws.Range('C29').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.160'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -7.19%  '

$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.09%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.63'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -8.98%  '

$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.84'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -10.16%  '

$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.84'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -11.19%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '21.41'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -5.37%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.19'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -6.87%  '

$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.76'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -9.39%  '

$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '155.80'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -6.31%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.10'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -9.86%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.37'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -11.96%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '25.25'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -5.62%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.149.86'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -5.23%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0678'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -8.95%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.15'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -4.46%  '

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -10.20%  '

$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.86'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -8.55%  '

$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.994'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.51%  '

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -8.57%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.45'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -9.59%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.228.33'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -5.40%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.08'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -7.27%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.81'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -7.14%  '
